$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 10
$ws.Range("F8").Value = -5
